$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cl = $m.CustomLayouts.Item(1)
$dt = $cl.HeadersFooters.DateAndTime
$dt.UseFormat = 0
$dt.Text = "27/09/2016"
$dt.Visible = -1
Write-Output "done"
